$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up row 19 (E19): drop the quote-prefix style, keep same text ---
# Copy the formatting used by the row above (E18, which already has the
# plain "vertical center + wrap" style) onto E19 so its stray
# quote-prefix/fill formatting goes away while its text is untouched.
$ws.Range("E18").Copy()
$ws.Range("E19").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Add the new meeting entry as row 20 ---
# Bring over the date/time/member formatting from row 19 first so the new
# row matches the existing table styling exactly.
$ws.Range("A19:E19").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(20, 1).Value = 45221
$ws.Cells.Item(20, 2).Value = 0.95833333333333337
$ws.Cells.Item(20, 3).Value = 0.99305555555555547
$ws.Cells.Item(20, 4).Value = "Aishwarya Anil Kumar (32644329) / Chih Hui Wang (33209006) / Shreyansh Mahtolia (33509115)"
$ws.Cells.Item(20, 5).Value = "'- Discuss and fix the error issues with the Chuck 19 left join`n- Discuss the next steps for the analysis section, next meeting, the date of completion and the date of recording`n- Discussion of the completed section: IDA and temporal analysis (keep or delete plots)"
$ws.Rows.Item(20).RowHeight = 85

# --- The old placeholder row 20 shifted down to row 21: it loses its spare C cell ---
$ws.Range("C21").Clear()

# --- Row 25 gains an (empty) formatted E cell matching the table's style ---
$ws.Range("E18").Copy()
$ws.Range("E25").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Update the view: scrolled down one row, new active selection ---
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("G20").Select() | Out-Null
